$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-02-09 Thursday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-02-10 Friday", 2) | Out-Null
$d.Content.Find.Execute("64-61=", $true, $true, $false, $false, $false, $true, 1, $false, "64-36=", 2) | Out-Null
$d.Content.Find.Execute("42-23=", $true, $true, $false, $false, $false, $true, 1, $false, "53-35=", 2) | Out-Null
$d.Content.Find.Execute("80-13=", $true, $true, $false, $false, $false, $true, 1, $false, "64-32=", 2) | Out-Null
$d.Content.Find.Execute("70-16=", $true, $true, $false, $false, $false, $true, 1, $false, "53-15=", 2) | Out-Null
$d.Content.Find.Execute("0+94=", $true, $true, $false, $false, $false, $true, 1, $false, "72-58=", 2) | Out-Null
$d.Content.Find.Execute("10+8=", $true, $true, $false, $false, $false, $true, 1, $false, "43+1=", 2) | Out-Null
$d.Content.Find.Execute("75-53=", $true, $true, $false, $false, $false, $true, 1, $false, "54-48=", 2) | Out-Null
$d.Content.Find.Execute("50+33=", $true, $true, $false, $false, $false, $true, 1, $false, "43-0=", 2) | Out-Null
$d.Content.Find.Execute("15+57=", $true, $true, $false, $false, $false, $true, 1, $false, "8+57=", 2) | Out-Null
$d.Content.Find.Execute("71+26=", $true, $true, $false, $false, $false, $true, 1, $false, "46-11=", 2) | Out-Null
$d.Content.Find.Execute("32+0=", $true, $true, $false, $false, $false, $true, 1, $false, "13-0=", 2) | Out-Null
$d.Content.Find.Execute("18+38=", $true, $true, $false, $false, $false, $true, 1, $false, "11+55=", 2) | Out-Null
$d.Content.Find.Execute("53+24=", $true, $true, $false, $false, $false, $true, 1, $false, "95-50=", 2) | Out-Null
$d.Content.Find.Execute("98-47=", $true, $true, $false, $false, $false, $true, 1, $false, "44+13=", 2) | Out-Null
$d.Content.Find.Execute("6+81=", $true, $true, $false, $false, $false, $true, 1, $false, "83-69=", 2) | Out-Null
$d.Content.Find.Execute("58+25=", $true, $true, $false, $false, $false, $true, 1, $false, "84-22=", 2) | Out-Null
$d.Content.Find.Execute("37+44=", $true, $true, $false, $false, $false, $true, 1, $false, "9+53=", 2) | Out-Null
$d.Content.Find.Execute("5+80=", $true, $true, $false, $false, $false, $true, 1, $false, "66-60=", 2) | Out-Null
$d.Content.Find.Execute("13+55=", $true, $true, $false, $false, $false, $true, 1, $false, "25+5=", 2) | Out-Null
$d.Content.Find.Execute("77-64=", $true, $true, $false, $false, $false, $true, 1, $false, "20+46=", 2) | Out-Null
$d.Content.Find.Execute("92-78=", $true, $true, $false, $false, $false, $true, 1, $false, "70-48=", 2) | Out-Null
$d.Content.Find.Execute("87-64=", $true, $true, $false, $false, $false, $true, 1, $false, "22+60=", 2) | Out-Null
$d.Content.Find.Execute("93-2=", $true, $true, $false, $false, $false, $true, 1, $false, "38-28=", 2) | Out-Null
$d.Content.Find.Execute("87-33=", $true, $true, $false, $false, $false, $true, 1, $false, "46+2=", 2) | Out-Null
$d.Content.Find.Execute("88-25=", $true, $true, $false, $false, $false, $true, 1, $false, "81-76=", 2) | Out-Null
$d.Content.Find.Execute("60-20=", $true, $true, $false, $false, $false, $true, 1, $false, "15+18=", 2) | Out-Null
$d.Content.Find.Execute("1+39=", $true, $true, $false, $false, $false, $true, 1, $false, "7+92=", 2) | Out-Null
$d.Content.Find.Execute("72-14=", $true, $true, $false, $false, $false, $true, 1, $false, "51-7=", 2) | Out-Null
$d.Content.Find.Execute("29-20=", $true, $true, $false, $false, $false, $true, 1, $false, "58-24=", 2) | Out-Null
$d.Content.Find.Execute("54+35=", $true, $true, $false, $false, $false, $true, 1, $false, "12+11=", 2) | Out-Null
$d.Content.Find.Execute("81-8=", $true, $true, $false, $false, $false, $true, 1, $false, "11+69=", 2) | Out-Null
$d.Content.Find.Execute("47+15=", $true, $true, $false, $false, $false, $true, 1, $false, "73-38=", 2) | Out-Null
$d.Content.Find.Execute("62+32=", $true, $true, $false, $false, $false, $true, 1, $false, "77+3=", 2) | Out-Null
$d.Content.Find.Execute("48-34=", $true, $true, $false, $false, $false, $true, 1, $false, "63-20=", 2) | Out-Null
$d.Content.Find.Execute("58-6=", $true, $true, $false, $false, $false, $true, 1, $false, "21-3=", 2) | Out-Null
$d.Content.Find.Execute("84+1=", $true, $true, $false, $false, $false, $true, 1, $false, "4+13=", 2) | Out-Null
$d.Content.Find.Execute("48-47=", $true, $true, $false, $false, $false, $true, 1, $false, "24+31=", 2) | Out-Null
$d.Content.Find.Execute("94-91=", $true, $true, $false, $false, $false, $true, 1, $false, "7+79=", 2) | Out-Null
$d.Content.Find.Execute("43-25=", $true, $true, $false, $false, $false, $true, 1, $false, "7+48=", 2) | Out-Null
$d.Content.Find.Execute("50+12=", $true, $true, $false, $false, $false, $true, 1, $false, "6+56=", 2) | Out-Null
$d.Content.Find.Execute("19+5=", $true, $true, $false, $false, $false, $true, 1, $false, "75+23=", 2) | Out-Null
$d.Content.Find.Execute("34+59=", $true, $true, $false, $false, $false, $true, 1, $false, "30+50=", 2) | Out-Null
$d.Content.Find.Execute("43+46=", $true, $true, $false, $false, $false, $true, 1, $false, "52+36=", 2) | Out-Null
$d.Content.Find.Execute("56+17=", $true, $true, $false, $false, $false, $true, 1, $false, "84+6=", 2) | Out-Null
$d.Content.Find.Execute("60+33=", $true, $true, $false, $false, $false, $true, 1, $false, "15+43=", 2) | Out-Null
$d.Content.Find.Execute("42+20=", $true, $true, $false, $false, $false, $true, 1, $false, "56-42=", 2) | Out-Null
$d.Content.Find.Execute("24+9=", $true, $true, $false, $false, $false, $true, 1, $false, "42+6=", 2) | Out-Null
$d.Content.Find.Execute("48-27=", $true, $true, $false, $false, $false, $true, 1, $false, "78-2=", 2) | Out-Null
$d.Content.Find.Execute("73+12=", $true, $true, $false, $false, $false, $true, 1, $false, "40-19=", 2) | Out-Null
$d.Content.Find.Execute("94-21=", $true, $true, $false, $false, $false, $true, 1, $false, "24-23=", 2) | Out-Null
$d.Content.Find.Execute("23-12=", $true, $true, $false, $false, $false, $true, 1, $false, "78+21=", 2) | Out-Null
$d.Content.Find.Execute("22+38=", $true, $true, $false, $false, $false, $true, 1, $false, "86-46=", 2) | Out-Null
$d.Content.Find.Execute("23+19=", $true, $true, $false, $false, $false, $true, 1, $false, "66-27=", 2) | Out-Null
$d.Content.Find.Execute("12+86=", $true, $true, $false, $false, $false, $true, 1, $false, "92-73=", 2) | Out-Null
$d.Content.Find.Execute("88-80=", $true, $true, $false, $false, $false, $true, 1, $false, "9-7=", 2) | Out-Null
$d.Content.Find.Execute("15+15=", $true, $true, $false, $false, $false, $true, 1, $false, "17+33=", 2) | Out-Null
$d.Content.Find.Execute("68+22=", $true, $true, $false, $false, $false, $true, 1, $false, "92-71=", 2) | Out-Null
$d.Content.Find.Execute("70-9=", $true, $true, $false, $false, $false, $true, 1, $false, "54-29=", 2) | Out-Null
$d.Content.Find.Execute("66+2=", $true, $true, $false, $false, $false, $true, 1, $false, "95-9=", 2) | Out-Null
$d.Content.Find.Execute("49-23=", $true, $true, $false, $false, $false, $true, 1, $false, "5+37=", 2) | Out-Null
$d.Content.Find.Execute("28+1=", $true, $true, $false, $false, $false, $true, 1, $false, "69-54=", 2) | Out-Null
$d.Content.Find.Execute("83+6=", $true, $true, $false, $false, $false, $true, 1, $false, "17+3=", 2) | Out-Null
$d.Content.Find.Execute("53+20=", $true, $true, $false, $false, $false, $true, 1, $false, "59-6=", 2) | Out-Null
$d.Content.Find.Execute("50+35=", $true, $true, $false, $false, $false, $true, 1, $false, "1+87=", 2) | Out-Null
$d.Content.Find.Execute("71-13=", $true, $true, $false, $false, $false, $true, 1, $false, "24+45=", 2) | Out-Null
$d.Content.Find.Execute("18+37=", $true, $true, $false, $false, $false, $true, 1, $false, "59+26=", 2) | Out-Null
$d.Content.Find.Execute("82-42=", $true, $true, $false, $false, $false, $true, 1, $false, "4+84=", 2) | Out-Null
$d.Content.Find.Execute("51-13=", $true, $true, $false, $false, $false, $true, 1, $false, "78-30=", 2) | Out-Null
$d.Content.Find.Execute("23-8=", $true, $true, $false, $false, $false, $true, 1, $false, "96-55=", 2) | Out-Null
$d.Content.Find.Execute("53-7=", $true, $true, $false, $false, $false, $true, 1, $false, "99-55=", 2) | Out-Null
$d.Content.Find.Execute("86-63=", $true, $true, $false, $false, $false, $true, 1, $false, "60-11=", 2) | Out-Null
$d.Content.Find.Execute("11+49=", $true, $true, $false, $false, $false, $true, 1, $false, "20-1=", 2) | Out-Null
$d.Content.Find.Execute("67-33=", $true, $true, $false, $false, $false, $true, 1, $false, "63+24=", 2) | Out-Null
$d.Content.Find.Execute("92-53=", $true, $true, $false, $false, $false, $true, 1, $false, "69+0=", 2) | Out-Null
$d.Content.Find.Execute("81+11=", $true, $true, $false, $false, $false, $true, 1, $false, "84-84=", 2) | Out-Null
$d.Content.Find.Execute("61+36=", $true, $true, $false, $false, $false, $true, 1, $false, "29+1=", 2) | Out-Null
$d.Content.Find.Execute("27-26=", $true, $true, $false, $false, $false, $true, 1, $false, "20+61=", 2) | Out-Null
$d.Content.Find.Execute("86-52=", $true, $true, $false, $false, $false, $true, 1, $false, "99-80=", 2) | Out-Null
$d.Content.Find.Execute("24+38=", $true, $true, $false, $false, $false, $true, 1, $false, "16+83=", 2) | Out-Null
$d.Content.Find.Execute("17+24=", $true, $true, $false, $false, $false, $true, 1, $false, "13+9=", 2) | Out-Null
$d.Content.Find.Execute("72-63=", $true, $true, $false, $false, $false, $true, 1, $false, "4-0=", 2) | Out-Null
$d.Content.Find.Execute("95-61=", $true, $true, $false, $false, $false, $true, 1, $false, "53-14=", 2) | Out-Null
$d.Content.Find.Execute("57+12=", $true, $true, $false, $false, $false, $true, 1, $false, "52-26=", 2) | Out-Null
$d.Content.Find.Execute("64+22=", $true, $true, $false, $false, $false, $true, 1, $false, "82-81=", 2) | Out-Null
$d.Content.Find.Execute("48+28=", $true, $true, $false, $false, $false, $true, 1, $false, "48+44=", 2) | Out-Null
$d.Content.Find.Execute("8+40=", $true, $true, $false, $false, $false, $true, 1, $false, "84-12=", 2) | Out-Null
$d.Content.Find.Execute("25+0=", $true, $true, $false, $false, $false, $true, 1, $false, "76-37=", 2) | Out-Null
$d.Content.Find.Execute("98-10=", $true, $true, $false, $false, $false, $true, 1, $false, "44+17=", 2) | Out-Null
$d.Content.Find.Execute("76-17=", $true, $true, $false, $false, $false, $true, 1, $false, "4+66=", 2) | Out-Null
$d.Content.Find.Execute("14+47=", $true, $true, $false, $false, $false, $true, 1, $false, "64+0=", 2) | Out-Null
$d.Content.Find.Execute("80-40=", $true, $true, $false, $false, $false, $true, 1, $false, "85-44=", 2) | Out-Null
$d.Content.Find.Execute("87-75=", $true, $true, $false, $false, $false, $true, 1, $false, "69-52=", 2) | Out-Null
$d.Content.Find.Execute("83+2=", $true, $true, $false, $false, $false, $true, 1, $false, "24-1=", 2) | Out-Null
$d.Content.Find.Execute("90-72=", $true, $true, $false, $false, $false, $true, 1, $false, "8+36=", 2) | Out-Null
$d.Content.Find.Execute("27+8=", $true, $true, $false, $false, $false, $true, 1, $false, "12+74=", 2) | Out-Null
$d.Content.Find.Execute("12+79=", $true, $true, $false, $false, $false, $true, 1, $false, "26-25=", 2) | Out-Null
$d.Content.Find.Execute("40-11=", $true, $true, $false, $false, $false, $true, 1, $false, "1+74=", 2) | Out-Null
$d.Content.Find.Execute("87-71=", $true, $true, $false, $false, $false, $true, 1, $false, "4+83=", 2) | Out-Null
$d.Content.Find.Execute("94-77=", $true, $true, $false, $false, $false, $true, 1, $false, "90-39=", 2) | Out-Null
$d.Content.Find.Execute("57+15=", $true, $true, $false, $false, $false, $true, 1, $false, "39+24=", 2) | Out-Null
